# Insert a new data row at row 247. This shifts the existing rows 247:311
# down to 248:312 (values, styles, formats all move with the rows), then
# populate the freshly-inserted row 247 with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("247:247").Insert()

$ws.Range("A247").Value = 8
$ws.Range("B247").Value = "Terminal La Palmera de La Serena"
$ws.Range("C247").Value = "Coquimbo"
$ws.Range("D247").Value = 44642
$ws.Range("E247").Value = 4
$ws.Range("F247").Value = 100114013
$ws.Range("G247").Value = "Zanahoria"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 500
$ws.Range("K247").Value = 5800
$ws.Range("L247").Value = 6000
$ws.Range("M247").Value = 5900
$ws.Range("N247").Value = "$/saco 20 kilos"
$ws.Range("O247").Value = "Provincia del Elquí"
$ws.Range("P247").Value = 295
$ws.Range("Q247").Value = 20
$ws.Range("R247").Value = "Hortaliza"

# Match the date-serial number format used by the other rows' Fecha column.
$ws.Range("D247").NumberFormat = $ws.Range("D248").NumberFormat
